$d = $word.ActiveDocument

$d.Content.Find.Execute("339÷7=", $true, $false, $false, $false, $false, $true, 1, $false, "963÷4=", 2) | Out-Null
$d.Content.Find.Execute("178÷9=", $true, $false, $false, $false, $false, $true, 1, $false, "628÷5=", 2) | Out-Null
$d.Content.Find.Execute("806÷6=", $true, $false, $false, $false, $false, $true, 1, $false, "207÷4=", 2) | Out-Null
$d.Content.Find.Execute("669÷5=", $true, $false, $false, $false, $false, $true, 1, $false, "296÷2=", 2) | Out-Null
$d.Content.Find.Execute("884÷2=", $true, $false, $false, $false, $false, $true, 1, $false, "158÷7=", 2) | Out-Null
$d.Content.Find.Execute("540÷5=", $true, $false, $false, $false, $false, $true, 1, $false, "800÷9=", 2) | Out-Null
$d.Content.Find.Execute("251÷4=", $true, $false, $false, $false, $false, $true, 1, $false, "434÷2=", 2) | Out-Null
$d.Content.Find.Execute("777÷5=", $true, $false, $false, $false, $false, $true, 1, $false, "460÷2=", 2) | Out-Null
$d.Content.Find.Execute("675÷6=", $true, $false, $false, $false, $false, $true, 1, $false, "992÷3=", 2) | Out-Null
$d.Content.Find.Execute("950÷6=", $true, $false, $false, $false, $false, $true, 1, $false, "849÷2=", 2) | Out-Null
$d.Content.Find.Execute("634÷5=", $true, $false, $false, $false, $false, $true, 1, $false, "755÷3=", 2) | Out-Null
$d.Content.Find.Execute("293÷6=", $true, $false, $false, $false, $false, $true, 1, $false, "862÷8=", 2) | Out-Null
$d.Content.Find.Execute("505÷5=", $true, $false, $false, $false, $false, $true, 1, $false, "940÷4=", 2) | Out-Null
$d.Content.Find.Execute("345÷5=", $true, $false, $false, $false, $false, $true, 1, $false, "138÷6=", 2) | Out-Null
$d.Content.Find.Execute("721÷3=", $true, $false, $false, $false, $false, $true, 1, $false, "205÷7=", 2) | Out-Null
$d.Content.Find.Execute("749÷5=", $true, $false, $false, $false, $false, $true, 1, $false, "986÷7=", 2) | Out-Null
$d.Content.Find.Execute("921÷8=", $true, $false, $false, $false, $false, $true, 1, $false, "806÷6=", 2) | Out-Null
$d.Content.Find.Execute("886÷8=", $true, $false, $false, $false, $false, $true, 1, $false, "695÷6=", 2) | Out-Null
$d.Content.Find.Execute("496÷9=", $true, $false, $false, $false, $false, $true, 1, $false, "111÷4=", 2) | Out-Null
$d.Content.Find.Execute("196÷9=", $true, $false, $false, $false, $false, $true, 1, $false, "315÷3=", 2) | Out-Null
$d.Content.Find.Execute("305÷8=", $true, $false, $false, $false, $false, $true, 1, $false, "722÷5=", 2) | Out-Null
$d.Content.Find.Execute("252÷7=", $true, $false, $false, $false, $false, $true, 1, $false, "847÷3=", 2) | Out-Null
$d.Content.Find.Execute("106÷5=", $true, $false, $false, $false, $false, $true, 1, $false, "881÷7=", 2) | Out-Null
$d.Content.Find.Execute("523÷6=", $true, $false, $false, $false, $false, $true, 1, $false, "639÷2=", 2) | Out-Null
$d.Content.Find.Execute("199÷6=", $true, $false, $false, $false, $false, $true, 1, $false, "851÷3=", 2) | Out-Null
